$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Style = $ws.Range("H1").Style

# Add data values for columns I and J, rows 2-24
$data = @(
    @(6, 7),
    @(8, 9),
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(5, 7),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(4, 5),
    @(8, 9),
    @(8, 8),
    @(4, 5),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(5, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
